# Add a new "brokersCanContinueRenewalSubmission" test-data block to the
# DashboardPageData sheet, mirroring the existing blocks (header row +
# label row + two data rows), as described by the commit
# "added test for renewal submission".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DashboardPageData")

# Clone the formatting of the most similar existing block (rows 55-58,
# "sortQuoteList") down onto the new rows 74-77, leaving a one-row gap
# (row 72/73) just like the other blocks on this sheet.
$ws.Range("A55:B55").Copy($ws.Range("A74:B74"))
$ws.Range("A56:D56").Copy($ws.Range("A75:D75"))
$ws.Range("A57:D57").Copy($ws.Range("A76:D76"))
$ws.Range("A58:D58").Copy($ws.Range("A77:D77"))

# New block header / label (this introduces the new shared string).
$ws.Range("A74").Value = "brokersCanContinueRenewalSubmission"

# Row 75: column headers for this block.
$ws.Range("A75").Value = "runMode"
$ws.Range("B75").Value = "brokerId"
$ws.Range("C75").Value = "agentId"
$ws.Range("D75").Value = "agencyOfficeId"

# Row 76: "Y" run-mode data row.
$ws.Range("A76").Value = "Y"
$ws.Range("B76").Value = 20217
$ws.Range("C76").Value = 237
$ws.Range("D76").Value = 8006

# Row 77: "N" run-mode data row.
$ws.Range("A77").Value = "N"
$ws.Range("B77").Value = 25997
$ws.Range("C77").Value = 7166
$ws.Range("D77").Value = 8414

# Reflect the new selection left on this sheet after entering the data.
$ws.Range("A74:E78").Select()
